# Update odds values in the "Jogos do Dia" worksheet to reflect the
# latest Betfair back/lay prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sydney FC vs Auckland FC)
$ws.Range("F2").Value = 2.82
$ws.Range("H2").Value = 2.54
$ws.Range("J2").Value = 3.6
$ws.Range("K2").Value = 3.85
$ws.Range("L2").Value = 1.36
$ws.Range("N2").Value = 4.4
$ws.Range("P2").Value = 2.18
$ws.Range("R2").Value = 1.47
$ws.Range("T2").Value = 1.65
$ws.Range("V2").Value = 1.6
$ws.Range("W2").Value = 1.5
$ws.Range("Y2").Value = 13.5
$ws.Range("Z2").Value = 22
$ws.Range("AB2").Value = 14.5
$ws.Range("AD2").Value = 12.5
$ws.Range("AE2").Value = 28
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 13.5
$ws.Range("AI2").Value = 1000
$ws.Range("AM2").Value = 70
$ws.Range("AO2").Value = 19.5

# Row 4
$ws.Range("H4").Value = 2.48

# Row 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 4.8
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 4.7
$ws.Range("P5").Value = 2.34
$ws.Range("Q5").Value = 1.6

# Row 6
$ws.Range("F6").Value = 3.7
$ws.Range("H6").Value = 2.16
$ws.Range("I6").Value = 2.18
$ws.Range("AJ6").Value = 80

# Row 7
$ws.Range("F7").Value = 2.26
$ws.Range("G7").Value = 2.56
$ws.Range("J7").Value = 3.3
$ws.Range("P7").Value = 1.89

# Row 8
$ws.Range("I8").Value = 3.8
$ws.Range("P8").Value = 1.96
$ws.Range("Q8").Value = 1.82

# Row 9
$ws.Range("F9").Value = 1.79
$ws.Range("G9").Value = 1.86
$ws.Range("H9").Value = 4.7
$ws.Range("I9").Value = 5.2
$ws.Range("P9").Value = 2.04
$ws.Range("Q9").Value = 1.67

# Row 10
$ws.Range("Q10").Value = 1.65

# Row 11
$ws.Range("Q11").Value = 1.65

# Row 12
$ws.Range("I12").Value = 3.9
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 3.9
$ws.Range("T12").Value = 1.65
$ws.Range("Y12").Value = 18
$ws.Range("AA12").Value = 70
$ws.Range("AC12").Value = 8.800000000000001
$ws.Range("AE12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AK12").Value = 20
$ws.Range("AL12").Value = 34
$ws.Range("AO12").Value = 36

# Row 13
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 6.6
$ws.Range("Z13").Value = 11.5

# Row 14
$ws.Range("U14").Value = 2.62
$ws.Range("X14").Value = 26
$ws.Range("AF14").Value = 34

# Row 16
$ws.Range("P16").Value = 1.61
